$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.809.16"
$ws.Range("E2").Value = "  +4.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.626.79"
$ws.Range("E3").Value = "  +3.19%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "629.67"
$ws.Range("E5").Value = "  +3.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.63"
$ws.Range("E6").Value = "  +5.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.623.77"
$ws.Range("E7").Value = "  +3.02%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.496"
$ws.Range("E9").Value = "  +3.07%  "
$ws.Range("E10").Value = "  +6.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.32"
$ws.Range("E11").Value = "  +6.16%  "
$ws.Range("E12").Value = "  +3.79%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000228"
$ws.Range("E13").Value = "  +4.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.45"
$ws.Range("E14").Value = "  +5.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.242.37"
$ws.Range("E15").Value = "  +3.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.645.42"
$ws.Range("E16").Value = "  +3.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.914.54"
$ws.Range("E17").Value = "  +4.55%  "
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.68"
$ws.Range("E19").Value = "  +6.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.04"
$ws.Range("E20").Value = "  +4.81%  "
$ws.Range("E21").Value = "  +11.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "464.89"
$ws.Range("E22").Value = "  +4.66%  "
$ws.Range("E23").Value = "  +3.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.00"
$ws.Range("E24").Value = "  +2.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000136"
$ws.Range("E25").Value = "  +12.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.77"
$ws.Range("E26").Value = "  +5.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.773.20"
$ws.Range("E27").Value = "  +3.31%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.29"
$ws.Range("E29").Value = "  +14.21%  "
$ws.Range("E30").Value = "  +4.93%  "
$ws.Range("E31").Value = "  +6.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.180"
$ws.Range("E32").Value = "  +12.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.63"
$ws.Range("E33").Value = "  +8.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("E35").Value = "  +5.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.65"
$ws.Range("E36").Value = "  +4.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.625.74"
$ws.Range("E37").Value = "  +3.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.53"
$ws.Range("E38").Value = "  +6.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.44"
$ws.Range("E39").Value = "  +13.76%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0936"
$ws.Range("E41").Value = "  +8.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "179.76"
$ws.Range("E42").Value = "  +3.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.69"
$ws.Range("E44").Value = "  +2.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "31.89"
$ws.Range("E45").Value = "  +17.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.915"
$ws.Range("E46").Value = "  +2.75%  "
$ws.Range("E47").Value = "  +12.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.79"
$ws.Range("E48").Value = "  +10.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "46.36"
$ws.Range("E49").Value = "  +2.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.83"
$ws.Range("E50").Value = "  +3.70%  "
$ws.Range("E51").Value = "  +9.35%  "
